$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("http://purl.obolibrary.org/obo/BFO_0000003", "occurrent", "R:min", "SubClassOf: ns1:BFO_0000001 | SubClassOf: ns1:BFO_0000132 min 0", "http://ies.data.gov.uk/ontology/ies4#Entity", "Entity", "R:min", "SubClassOf: ies:Element | SubClassOf: ies:isPartOf min 0"),
    @("http://purl.obolibrary.org/obo/BFO_0000003", "occurrent", "R:min", "SubClassOf: ns1:BFO_0000001 | SubClassOf: ns1:BFO_0000132 min 0", "http://ies.data.gov.uk/ontology/ies4#State", "State", "R:min", "SubClassOf: ies:isStateOf min 1"),
    @("http://purl.obolibrary.org/obo/BFO_0000003", "occurrent", "R:min", "SubClassOf: ns1:BFO_0000001 | SubClassOf: ns1:BFO_0000132 min 0", "http://ies.data.gov.uk/ontology/ies4#Event", "Event", "R:min", "SubClassOf: ies:Element | SubClassOf: ies:isPartOf min 0")
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
